# Integrate Dwarven light armor
# Insert a new "Light_Dwarven" row into the Armors sheet, just above the
# existing "Light_Falmer" row (current row 70), shifting every following
# row down by one, then refresh the sheet's bookkeeping (dimension is
# automatic, selection + sort range need to be restated).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Armors")

# --- Insert the new row ------------------------------------------------
$ws.Rows.Item(70).Insert()

# The insert copies formatting down from the row above; the new armor
# row should be unstyled (like its neighbours), so strip that back off.
$ws.Range("A70:I70").ClearFormats()

# --- Populate the new row ----------------------------------------------
$ws.Range("A70").Value = "Light_Dwarven"
$ws.Range("B70").Value = 400
$ws.Range("C70").Value = 30
$ws.Range("D70").Value = 1000
$ws.Range("E70").Value = "Dwarven"
$ws.Range("F70").Value = "Leather"
$ws.Range("G70").Value = "Leather Strips"
# No Temper (column H) for this entry.
$ws.Range("H70").ClearContents()
$ws.Range("I70").Value = "Dwarven Smithing"

# --- Restate the sort range / sort state so it covers the new row ------
$sortObj = $ws.Sort()
$sortObj.SortFields().Clear()
$sortObj.SortFields().Add($ws.Range("B42:B84"))
$sortObj.SetRange($ws.Range("A42:I84"))
$sortObj.Header = 1
$sortObj.Apply()

# --- Restore selection to match the author's final cursor position -----
$ws.Activate()
[void]$ws.Range("J69").Select()
